$wb = $excel.ActiveWorkbook

# ALC row 15 (G15=44146)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2757.7778
$ws.Range("I15").Value = 2757.7778
$ws.Range("K15").Value = 8273.3334
$ws.Range("M15").Value = -8104.3334

# ALC row 55 (G55=5517)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 167
$ws.Range("I55").Value = 93.333336
$ws.Range("J55").Value = 255.4
$ws.Range("K55").Value = 93.333336
$ws.Range("L55").Value = 255.4
$ws.Range("M55").Value = 120.666664
$ws.Range("N55").Value = -683.4

# ALC row 86 (G86=12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3170.8333
$ws.Range("I86").Value = 2561.3635
$ws.Range("K86").Value = 2561.3635
$ws.Range("M86").Value = -1438.3635

# ALC row 89 (G89=12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3170.8333
$ws.Range("I89").Value = 2561.3635
$ws.Range("K89").Value = 12806.8175
$ws.Range("M89").Value = -7190.817499999999

# ALC row 113 (G113=27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 11076.923
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 11750
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 11750
$ws.Range("M113").Value = -6746
$ws.Range("N113").Value = -18258

# ALC row 116 (G116=27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9232.137000000001
$ws.Range("I116").Value = 8607.111000000001
$ws.Range("J116").Value = 9664.846
$ws.Range("K116").Value = 8607.111000000001
$ws.Range("L116").Value = 9664.846
$ws.Range("M116").Value = -5165.111000000001
$ws.Range("N116").Value = -16548.846

# ALC row 123 (G123=34090)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 54454.547
$ws.Range("J123").Value = 54454.547
$ws.Range("L123").Value = 54454.547
$ws.Range("N123").Value = -64254.547

# ALC row 124 (G124=34241)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820

# ALC row 128 (G128=34540)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

# ALC row 132 (G132=44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 15496.609
$ws.Range("I132").Value = 3141.423
$ws.Range("J132").Value = 36912.266
$ws.Range("K132").Value = 9424.269
$ws.Range("L132").Value = 110736.798
$ws.Range("M132").Value = -6894.269
$ws.Range("N132").Value = -115796.798

# ALC row 136 (G136=42164)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 40000
$ws.Range("N136").Value = -50200

# ARM row 61 (G61=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1835.8889
$ws.Range("I61").Value = 1835.8889
$ws.Range("K61").Value = 1835.8889
$ws.Range("M61").Value = -1623.8889

# ARM row 63 (G63=12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6724.8
$ws.Range("I63").Value = 5749.8335
$ws.Range("J63").Value = 8187.25
$ws.Range("K63").Value = 5749.8335
$ws.Range("L63").Value = 8187.25
$ws.Range("M63").Value = -5063.8335
$ws.Range("N63").Value = -9559.25

# ARM row 66 (G66=12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6724.8
$ws.Range("I66").Value = 5749.8335
$ws.Range("J66").Value = 8187.25
$ws.Range("K66").Value = 28749.1675
$ws.Range("L66").Value = 40936.25
$ws.Range("M66").Value = -25317.1675
$ws.Range("N66").Value = -47800.25

# ARM row 97 (G97=19941)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 15205.6
$ws.Range("I97").Value = 14007
$ws.Range("K97").Value = 14007
$ws.Range("M97").Value = -13511

# ARM row 132 (G132=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2874.5
$ws.Range("I132").Value = 2874.5
$ws.Range("K132").Value = 8623.5
$ws.Range("M132").Value = -6093.5

# ARM row 136 (G136=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1835.8889
$ws.Range("I136").Value = 1835.8889
$ws.Range("K136").Value = 5507.6667
$ws.Range("M136").Value = -2957.6667

# BSM row 15 (G15=1605)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# BSM row 107 (G107=27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 9787.951999999999
$ws.Range("I107").Value = 10406.786
$ws.Range("K107").Value = 10406.786
$ws.Range("M107").Value = -8486.786

# CRP row 31 (G31=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1932.5294
$ws.Range("J31").Value = 1966.3334
$ws.Range("L31").Value = 1966.3334
$ws.Range("N31").Value = -2556.3334

# CRP row 34 (G34=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1932.5294
$ws.Range("J34").Value = 1966.3334
$ws.Range("L34").Value = 1966.3334
$ws.Range("N34").Value = -2370.3334

# CRP row 58 (G58=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2436.7856
$ws.Range("I58").Value = 1968.2
$ws.Range("K58").Value = 1968.2
$ws.Range("M58").Value = -1765.2

# CRP row 120 (G120=27230)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 79666.664
$ws.Range("J120").Value = 79666.664
$ws.Range("L120").Value = 79666.664
$ws.Range("N120").Value = -86924.664

# CRP row 122 (G122=36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3779
$ws.Range("I122").Value = 3634.7273
$ws.Range("J122").Value = 3862.5264
$ws.Range("K122").Value = 10904.1819
$ws.Range("L122").Value = 11587.5792
$ws.Range("M122").Value = -8454.1819
$ws.Range("N122").Value = -16487.5792

# CRP row 136 (G136=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2436.7856
$ws.Range("I136").Value = 1968.2
$ws.Range("K136").Value = 5904.6
$ws.Range("M136").Value = -3354.6

# CUL row 107 (G107=27838)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3311.75
$ws.Range("J107").Value = 3498.8
$ws.Range("L107").Value = 10496.4
$ws.Range("N107").Value = -14336.4

# CUL row 132 (G132=43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5650
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 1300
$ws.Range("K132").Value = 90000
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -87470
$ws.Range("N132").Value = -16760

# GSM row 107 (G107=27802)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 530.3889
$ws.Range("I107").Value = 400.23077
$ws.Range("K107").Value = 400.23077
$ws.Range("M107").Value = 1519.76923

# GSM row 113 (G113=27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 17499
$ws.Range("I113").Value = 14999.5
$ws.Range("J113").Value = 18748.75
$ws.Range("K113").Value = 14999.5
$ws.Range("L113").Value = 18748.75
$ws.Range("M113").Value = -12829.5
$ws.Range("N113").Value = -23088.75

# GSM row 122 (G122=36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3353
$ws.Range("I122").Value = 1039.625
$ws.Range("J122").Value = 6437.5
$ws.Range("K122").Value = 3118.875
$ws.Range("L122").Value = 19312.5
$ws.Range("M122").Value = -668.875
$ws.Range("N122").Value = -24212.5

# GSM row 126 (G126=36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5503.6665
$ws.Range("I126").Value = 5503.6665
$ws.Range("K126").Value = 16510.9995
$ws.Range("M126").Value = -14040.9995

# LTW row 16 (G16=5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 22727800
$ws.Range("I16").Value = 31250378
$ws.Range("K16").Value = 31250378
$ws.Range("M16").Value = -31250208

# LTW row 40 (G40=36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3163.4
$ws.Range("J40").Value = 3702.5
$ws.Range("L40").Value = 3702.5
$ws.Range("N40").Value = -3974.5

# LTW row 46 (G46=5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2715.08
$ws.Range("I46").Value = 1701.8182
$ws.Range("J46").Value = 3511.2144
$ws.Range("K46").Value = 1701.8182
$ws.Range("L46").Value = 3511.2144
$ws.Range("M46").Value = -1513.8182
$ws.Range("N46").Value = -3887.2144

# LTW row 93 (G93=19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2644.55
$ws.Range("I93").Value = 2395.5715
$ws.Range("J93").Value = 3225.5
$ws.Range("K93").Value = 2395.5715
$ws.Range("L93").Value = 3225.5
$ws.Range("M93").Value = -1147.5715
$ws.Range("N93").Value = -5721.5

# LTW row 122 (G122=36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5499.2
$ws.Range("I122").Value = 4165.3335
$ws.Range("K122").Value = 12496.0005
$ws.Range("M122").Value = -10046.0005

# WVR row 113 (G113=27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 461
$ws.Range("I113").Value = 439.2857
$ws.Range("J113").Value = 537
$ws.Range("K113").Value = 1317.8571
$ws.Range("L113").Value = 1611
$ws.Range("M113").Value = 852.1428999999998
$ws.Range("N113").Value = -5951

# WVR row 122 (G122=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1573.25
$ws.Range("I122").Value = 1599.3334
$ws.Range("J122").Value = 1495
$ws.Range("K122").Value = 4798.0002
$ws.Range("L122").Value = 4485
$ws.Range("M122").Value = -2348.0002
$ws.Range("N122").Value = -9385

# WVR row 126 (G126=36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6079.9414
$ws.Range("I126").Value = 6573.846
$ws.Range("K126").Value = 19721.538
$ws.Range("M126").Value = -17251.538

# WVR row 132 (G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3134.3572
$ws.Range("I132").Value = 2168.1
$ws.Range("K132").Value = 6504.299999999999
$ws.Range("M132").Value = -3974.299999999999
